$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy format from existing header cell (E1) then set value/text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F10 - time_taken values (no special style, matches rest of data rows)
$times = @(
    "2021-10-05 13:39:45.080059",
    "2021-10-05 13:39:45.080071",
    "2021-10-05 13:39:45.080075",
    "2021-10-05 13:39:45.080078",
    "2021-10-05 13:39:45.080082",
    "2021-10-05 13:39:45.080085",
    "2021-10-05 13:39:45.080088",
    "2021-10-05 13:39:45.080091",
    "2021-10-05 13:39:45.080094"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
